$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (9 rows: 3 sending clusters x 3 target clusters), replacing the
# previous 6 rows (2 sending clusters x 3 target clusters). Values recomputed
# with the 3rd cluster (ECs) now included in the analysis, per Dr Hou's advice.
$rows = @(
    ,@("ECs", "Tgfb3", "Eng", "ECs", 2, 0.6666666666666666, 0.6423476666666667, 1.927043, 0.01173234890143342, 0.01173234890143342, 3, 1, 114.155417, 342.466251, 0.6835107367845005, 0.6835107367845005, 73.32746574731033, 659.947191725793, 0.008019186441831581, 0.008019186441831583)
    ,@("ECs", "Tgfb3", "Eng", "FAPs", 2, 0.6666666666666666, 0.6423476666666667, 1.927043, 0.01173234890143342, 0.01173234890143342, 3, 1, 35.924535, 107.773605, 0.2150997826628812, 0.2150997826628812, 23.076041233335, 207.684371100015, 0.002523625698823421, 0.002523625698823422)
    ,@("ECs", "Tgfb3", "Eng", "sCs", 2, 0.6666666666666666, 0.6423476666666667, 1.927043, 0.01173234890143342, 0.01173234890143342, 3, 1, 16.93339666666667, 50.80019, 0.1013894805526183, 0.1013894805526183, 10.87712783757445, 97.89415053817001, 0.001189536760778417, 0.001189536760778417)
    ,@("FAPs", "Tgfb3", "Eng", "ECs", 3, 1, 22.27635266666667, 66.829058, 0.4068730304461968, 0.4068730304461968, 3, 1, 114.155417, 342.466251, 0.6835107367845005, 0.6835107367845005, 2542.966327902395, 22886.69695112156, 0.2781020848180225, 0.2781020848180225)
    ,@("FAPs", "Tgfb3", "Eng", "FAPs", 3, 1, 22.27635266666667, 66.829058, 0.4068730304461968, 0.4068730304461968, 3, 1, 35.924535, 107.773605, 0.2150997826628812, 0.2150997826628812, 800.26761104601, 7202.40849941409, 0.08751830042036476, 0.08751830042036478)
    ,@("FAPs", "Tgfb3", "Eng", "sCs", 3, 1, 22.27635266666667, 66.829058, 0.4068730304461968, 0.4068730304461968, 3, 1, 16.93339666666667, 50.80019, 0.1013894805526183, 0.1013894805526183, 377.2143159912245, 3394.92884392102, 0.04125264520780956, 0.04125264520780956)
    ,@("sCs", "Tgfb3", "Eng", "ECs", 3, 1, 31.831433, 95.494299, 0.5813946206523698, 0.5813946206523697, 3, 1, 114.155417, 342.466251, 0.6835107367845005, 0.6835107367845005, 3633.730507822561, 32703.57457040305, 0.3973894655246464, 0.3973894655246464)
    ,@("sCs", "Tgfb3", "Eng", "FAPs", 3, 1, 31.831433, 95.494299, 0.5813946206523698, 0.5813946206523697, 3, 1, 35.924535, 107.773605, 0.2150997826628812, 0.2150997826628812, 1143.529428908655, 10291.76486017789, 0.125057856543693, 0.125057856543693)
    ,@("sCs", "Tgfb3", "Eng", "sCs", 3, 1, 31.831433, 95.494299, 0.5813946206523698, 0.5813946206523697, 3, 1, 16.93339666666667, 50.80019, 0.1013894805526183, 0.1013894805526183, 539.0142814574233, 4851.12853311681, 0.05894729858403037, 0.05894729858403035)
)

# Write the recomputed LR-pair table back into the sheet starting at row 2
# (row 1 header is untouched). This overwrites the old 6 rows in place and
# extends the sheet down to row 10.
$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
